$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price-column cells we are about to rewrite so that
# numeric-looking strings (e.g. "202.50", "1.00") stay text instead of being
# auto-coerced into numbers (which would silently drop trailing zeros).
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D14","D15","D16","D17","D18","D19","D21","D22","D23","D24","D25","D26","D27","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46","D48","D49","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.499.81"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "3.603.99"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "202.50"
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("D6").Value = "595.63"
$ws.Range("E6").Value = "  -1.50%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("D9").Value = "0.215"
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("D10").Value = "0.644"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").Value = "53.82"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("B14").Value = "BitcoinCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D14").Value = "685.16"
$ws.Range("E14").Value = "  +15.90%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.173.33"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "70.563.14"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "19.10"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.610.84"
$ws.Range("E18").Value = "  +2.28%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").Value = "18.37"
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("D23").Value = "110.20"
$ws.Range("E23").Value = "  +6.85%  "
$ws.Range("D24").Value = "5.24"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "4.54"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "3.01"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").Value = "10.58"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "10.04"
$ws.Range("E29").Value = "  +5.75%  "
$ws.Range("D30").Value = "34.35"
$ws.Range("E30").Value = "  +3.84%  "
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("D32").Value = "7.15"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").Value = "12.32"
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").Value = "0.114"
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "63.66"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.861.98"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0848"
$ws.Range("E37").Value = "  +3.96%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "510.31"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").Value = "  -6.61%  "
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").Value = "36.72"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").Value = "0.383"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").Value = "0.137"
$ws.Range("D45").Value = "0.0470"
$ws.Range("E45").Value = "  +5.43%  "
$ws.Range("D46").Value = "3.06"
$ws.Range("E46").Value = "  +9.39%  "
$ws.Range("E47").Value = "  +3.63%  "
$ws.Range("D48").Value = "0.140"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").Value = "8.63"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "1.82"
$ws.Range("E51").Value = "  +23.63%  "
